# Update Name of Algo
# Applies updated numeric results to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.5381
$ws.Range("C3").Value = -11.5527
$ws.Range("C14").Value = -14.074
$ws.Range("C16").Value = -13.7103
$ws.Range("E18").Value = 18.09320000000001
$ws.Range("C21").Value = -12.65300000000001
$ws.Range("C23").Value = -12.73400000000001
$ws.Range("E24").Value = 16.4777
$ws.Range("C25").Value = -11.5923
$ws.Range("E25").Value = 16.9697
$ws.Range("C26").Value = -12.45410000000001
$ws.Range("E27").Value = 16.78519999999999
$ws.Range("C29").Value = -10.93080000000001
$ws.Range("E30").Value = 16.0967
$ws.Range("E31").Value = 16.6107
$ws.Range("E39").Value = 16.1339
$ws.Range("C40").Value = -12.4321
$ws.Range("E42").Value = 16.359
$ws.Range("E48").Value = 17.3917
$ws.Range("E51").Value = 17.0861
$ws.Range("E52").Value = 16.87210000000001
$ws.Range("C53").Value = -10.32020000000001
$ws.Range("E55").Value = 16.56050000000001
$ws.Range("E56").Value = 16.4041
$ws.Range("C57").Value = -14.10659999999999
$ws.Range("E57").Value = 16.911
$ws.Range("C59").Value = -12.6998
$ws.Range("E60").Value = 16.00440000000001
$ws.Range("C65").Value = -12.7495
$ws.Range("C69").Value = -11.20979999999999
$ws.Range("E73").Value = 17.49030000000001
$ws.Range("E74").Value = 16.95689999999998
$ws.Range("C79").Value = -11.38600000000001
$ws.Range("C83").Value = -13.9005
$ws.Range("E89").Value = 17.29520000000002
$ws.Range("E90").Value = 16.51529999999999
$ws.Range("C91").Value = -10.3218
$ws.Range("E92").Value = 18.79640000000002
$ws.Range("C93").Value = -12.1156
$ws.Range("C100").Value = -13.4225
